$d = $word.ActiveDocument

# 1. Merge "Product owner" text that was split across runs with spell-check
#    proofErr wrappers into a single run (text content is unchanged).
$d.Content.Find.Execute("Product owner", $false, $false, $false, $false, $false,
                         $true, 1, $false, "Product owner", 2) | Out-Null

# 2. Clear the paragraph containing the resolved doubt, leaving it empty.
$d.Content.Find.Execute("Levantamos a dúvida xxx que deverá ser esclarecida pelo professor Brandão.... ",
                         $false, $false, $false, $false, $false,
                         $true, 1, $false, "", 2) | Out-Null

# 3. Fill in the "Prazo" (deadline) column with "indefinido" for the four
#    populated task rows (rows 2-5) of the plan table, leaving the blank
#    trailing rows untouched.
$table = $d.Tables.Item(1)
for ($i = 2; $i -le 5; $i++) {
    $cell = $table.Cell($i, 3)
    $cell.Range.Text = "indefinido"
}
